# Actualiza dashboard y graficos
# Quita el sufijo de etapa (_ELE, _ENI, _ENIA, _EOC) de los nombres de los
# responsables en la columna C de la tabla de auditoria.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = "CamiloR"
    "C3"  = "CatalinaJ"
    "C4"  = "HoracioM"
    "C5"  = "RodolfoV"
    "C6"  = "ThiareM"
    "C7"  = "CatalinaJ"
    "C8"  = "HoracioM"
    "C9"  = "MarcelaP"
    "C10" = "RodolfoV"
    "C11" = "ThiareM"
    "C12" = "CatalinaJ"
    "C13" = "HoracioM"
    "C14" = "RodolfoV"
    "C15" = "ThiareM"
    "C16" = "CatalinaJ"
    "C17" = "HoracioM"
    "C18" = "CristianG"
    "C19" = "DanielaV"
    "C20" = "MarcelaP"
    "C21" = "CristianG"
    "C22" = "DanielaV"
    "C23" = "MarcelaP"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
